# Apply BOM edits: change LDO part (row 15) and correct its enable-pin
# designator, plus fix the 0R0 resistor designator on row 11.
# Also restore the view's zoom level and selected cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: 0R0 0603WAF0000T5E -> designator corrected from R7 to R1
$ws.Range("B11").Value = "R1"

# Row 15: LDO changed from AP7361C-33SP-13 (SO8-EP, U2) to AP7361C-33E-13
# (SOT-223, Q3) with corrected JLCPCB part number
$ws.Range("B15").Value = "Q3"
$ws.Range("A15").Value = "AP7361C-33E-13"
$ws.Range("C15").Value = "SOT-223"
$ws.Range("D15").Value = "C500795"

# View changes: zoom 115 -> 130, selection moved to C18
$ws.Application.ActiveWindow.Zoom = 130
$ws.Range("C18").Select()
